$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 9317
$ws.Range("E2").Value = 409
$ws.Range("F2").Value = 417
$ws.Range("G2").Value = 279
$ws.Range("H2").Value = 205
$ws.Range("I2").Value = 136
$ws.Range("J2").Value = 69
$ws.Range("K2").Value = 11184
$ws.Range("L2").Value = 6788
$ws.Range("M2").Value = 4397
$ws.Range("N2").Value = 2644
$ws.Range("O2").Value = 1753
$ws.Range("P2").Value = 499
$ws.Range("Q2").Value = 527
$ws.Range("R2").Value = -214
$ws.Range("S2").Value = -139
$ws.Range("T2").Value = 193
$ws.Range("U2").Value = 333
$ws.Range("V2").Value = 3816
$ws.Range("W2").Value = 4.39
$ws.Range("X2").Value = 2.2
$ws.Range("Y2").Value = 5.45
$ws.Range("Z2").Value = 1.84
$ws.Range("AA2").Value = 154.38
$ws.Range("AB2").Value = 469.24
$ws.Range("AC2").Value = 275
$ws.Range("AD2").Value = 24.54
$ws.Range("AE2").Value = 5749
$ws.Range("AF2").Value = 1.18
$ws.Range("AG2").Value = 70
$ws.Range("AH2").Value = 1.04
$ws.Range("AI2").Value = 23.68
$ws.Range("AJ2").Value = 49347483
$ws.Range("D3").Value = 8909
$ws.Range("E3").Value = 300
$ws.Range("F3").Value = 306
$ws.Range("G3").Value = 231
$ws.Range("H3").Value = 210
$ws.Range("I3").Value = 128
$ws.Range("J3").Value = 82
$ws.Range("K3").Value = 10565
$ws.Range("L3").Value = 5904
$ws.Range("M3").Value = 4661
$ws.Range("N3").Value = 2969
$ws.Range("O3").Value = 1691
$ws.Range("P3").Value = 499
$ws.Range("Q3").Value = 343
$ws.Range("R3").Value = -45
$ws.Range("S3").Value = -470
$ws.Range("T3").Value = 96
$ws.Range("U3").Value = 246
$ws.Range("V3").Value = 3399
$ws.Range("W3").Value = 3.37
$ws.Range("X3").Value = 2.36
$ws.Range("Y3").Value = 4.57
$ws.Range("Z3").Value = 1.93
$ws.Range("AA3").Value = 126.68
$ws.Range("AB3").Value = 531.36
$ws.Range("AC3").Value = 260
$ws.Range("AD3").Value = 45.98
$ws.Range("AE3").Value = 6457
$ws.Range("AF3").Value = 1.85
$ws.Range("AG3").Value = 70
$ws.Range("AH3").Value = 0.59
$ws.Range("AI3").Value = 25.1
$ws.Range("AJ3").Value = 49347483
$ws.Range("D4").Value = 8493
$ws.Range("E4").Value = 26
$ws.Range("F4").Value = 16
$ws.Range("G4").Value = -36
$ws.Range("H4").Value = -114
$ws.Range("I4").Value = -136
$ws.Range("J4").Value = 22
$ws.Range("K4").Value = 10658
$ws.Range("L4").Value = 6136
$ws.Range("M4").Value = 4522
$ws.Range("N4").Value = 2809
$ws.Range("O4").Value = 1714
$ws.Range("P4").Value = 499
$ws.Range("Q4").Value = -7
$ws.Range("R4").Value = 108
$ws.Range("S4").Value = 159
$ws.Range("T4").Value = 122
$ws.Range("U4").Value = -129
$ws.Range("V4").Value = 3637
$ws.Range("W4").Value = 0.31
$ws.Range("X4").Value = -1.35
$ws.Range("Y4").Value = -4.72
$ws.Range("Z4").Value = -1.08
$ws.Range("AA4").Value = 135.67
$ws.Range("AB4").Value = 498.56
$ws.Range("AC4").Value = -276
$ws.Range("AD4").Value = -23.39
$ws.Range("AE4").Value = 6108
$ws.Range("AF4").Value = 1.06
$ws.Range("AG4").Value = 30
$ws.Range("AH4").Value = 0.46
$ws.Range("AI4").Value = -10.12
$ws.Range("AJ4").Value = 49347483
$ws.Range("D5").Value = 9400
$ws.Range("E5").Value = 95
$ws.Range("F5").Value = 97
$ws.Range("G5").Value = -70
$ws.Range("H5").Value = -198
$ws.Range("I5").Value = -196
$ws.Range("J5").Value = -2
$ws.Range("K5").Value = 10793
$ws.Range("L5").Value = 6292
$ws.Range("M5").Value = 4502
$ws.Range("N5").Value = 2684
$ws.Range("O5").Value = 1817
$ws.Range("P5").Value = 499
$ws.Range("Q5").Value = 252
$ws.Range("R5").Value = -533
$ws.Range("S5").Value = 149
$ws.Range("T5").Value = 300
$ws.Range("U5").Value = -48
$ws.Range("V5").Value = 3615
$ws.Range("W5").Value = 1.01
$ws.Range("X5").Value = -2.11
$ws.Range("Y5").Value = -7.14
$ws.Range("Z5").Value = -1.85
$ws.Range("AA5").Value = 139.76
$ws.Range("AB5").Value = 478.43
$ws.Range("AC5").Value = -398
$ws.Range("AD5").Value = -13.15
$ws.Range("AE5").Value = 5837
$ws.Range("AF5").Value = 0.9
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 49347483
$ws.Range("D6").Value = 9294
$ws.Range("E6").Value = 235
$ws.Range("F6").Value = 241
$ws.Range("G6").Value = 104
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = -14
$ws.Range("K6").Value = 11166
$ws.Range("L6").Value = 6691
$ws.Range("M6").Value = 4475
$ws.Range("N6").Value = 2710
$ws.Range("P6").Value = 499
$ws.Range("Q6").Value = 68
$ws.Range("R6").Value = -197
$ws.Range("S6").Value = 377
$ws.Range("T6").Value = 262
$ws.Range("U6").Value = -195
$ws.Range("V6").Value = 4031
$ws.Range("W6").Value = 2.53
$ws.Range("X6").Value = 0.1
$ws.Range("Y6").Value = -0.53
$ws.Range("Z6").Value = 0.08
$ws.Range("AA6").Value = 149.5
$ws.Range("AB6").Value = 486.2
$ws.Range("AC6").Value = -29
$ws.Range("AD6").Value = -103.9
$ws.Range("AE6").Value = 5894
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 70
$ws.Range("AH6").Value = 2.31
$ws.Range("AI6").Value = -223.68
$ws.Range("AJ6").Value = 49347483
$ws.Range("D7").Value = 8140
$ws.Range("E7").Value = 360
$ws.Range("G7").Value = 370
$ws.Range("H7").Value = 290
$ws.Range("I7").Value = 180
$ws.Range("K7").Value = 11690
$ws.Range("L7").Value = 6450
$ws.Range("M7").Value = 5240
$ws.Range("N7").Value = 2920
$ws.Range("P7").Value = 500
$ws.Range("Q7").Value = 180
$ws.Range("R7").Value = -130
$ws.Range("S7").Value = 300
$ws.Range("T7").Value = 280
$ws.Range("W7").Value = 4.42
$ws.Range("X7").Value = 3.56
$ws.Range("Y7").Value = 6.39
$ws.Range("Z7").Value = 2.54
$ws.Range("AA7").Value = 123.09
$ws.Range("AC7").Value = 365
$ws.Range("AD7").Value = 12.57
$ws.Range("AE7").Value = 6350
$ws.Range("AF7").Value = 0.72
$ws.Range("AG7").Value = 70
$ws.Range("AH7").Value = 1.53
$ws.Range("AI7").Value = 19.19
$ws.Range("D8").Value = 8390
$ws.Range("E8").Value = 400
$ws.Range("G8").Value = 400
$ws.Range("H8").Value = 320
$ws.Range("I8").Value = 200
$ws.Range("K8").Value = 12040
$ws.Range("L8").Value = 6510
$ws.Range("M8").Value = 5530
$ws.Range("N8").Value = 3080
$ws.Range("P8").Value = 500
$ws.Range("Q8").Value = 390
$ws.Range("R8").Value = -160
$ws.Range("S8").Value = 30
$ws.Range("T8").Value = 270
$ws.Range("W8").Value = 4.77
$ws.Range("X8").Value = 3.81
$ws.Range("Y8").Value = 6.67
$ws.Range("Z8").Value = 2.7
$ws.Range("AA8").Value = 117.72
$ws.Range("AC8").Value = 405
$ws.Range("AD8").Value = 11.31
$ws.Range("AE8").Value = 6698
$ws.Range("AF8").Value = 0.68
$ws.Range("AG8").Value = 70
$ws.Range("AH8").Value = 1.53
$ws.Range("AI8").Value = 17.27
$ws.Range("D9").Value = 8810
$ws.Range("E9").Value = 460
$ws.Range("G9").Value = 440
$ws.Range("H9").Value = 350
$ws.Range("I9").Value = 220
$ws.Range("K9").Value = 12420
$ws.Range("L9").Value = 6580
$ws.Range("M9").Value = 5850
$ws.Range("N9").Value = 3270
$ws.Range("P9").Value = 500
$ws.Range("Q9").Value = 440
$ws.Range("R9").Value = -260
$ws.Range("S9").Value = 20
$ws.Range("T9").Value = 270
$ws.Range("W9").Value = 5.22
$ws.Range("X9").Value = 3.97
$ws.Range("Y9").Value = 6.93
$ws.Range("Z9").Value = 2.86
$ws.Range("AA9").Value = 112.48
$ws.Range("AC9").Value = 446
$ws.Range("AD9").Value = 10.28
$ws.Range("AE9").Value = 7111
$ws.Range("AF9").Value = 0.64
$ws.Range("AG9").Value = 70
$ws.Range("AH9").Value = 1.53
$ws.Range("AI9").Value = 15.7

$ws.Range("U7").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("U9").ClearContents()

Write-Host "Applied all changes"